{"js": "// Insert \", en total hay 3 sorteos diarios (esto puede variar)\" right\n// before the final period of the sentence that ends\n// \"... al igual que los clientes que ganaron en dicho sorteo.\"\nconst body = context.document.body;\n\n// Search for the sentence WITHOUT the trailing period so the collapsed\n// \"End\" caret lands exactly between \"sorteo\" and the final \".\".\nconst needle = \"al igual que los clientes que ganaron en dicho sorteo\";\nconst insertion = \", en total hay 3 sorteos diarios (esto puede variar)\";\n\nconst results = body.search(needle, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found in document body.\");\n}\n\nconst hit = results.items[0];\nconst caret = hit.getRange(\"End\");\ncaret.insertText(insertion, \"Before\");\n\nawait context.sync();\n", "ps1": "# Insert \", en total hay 3 sorteos diarios (esto puede variar)\" right\n# before the final period of the sentence that ends\n# \"... al igual que los clientes que ganaron en dicho sorteo.\"\n$d = $word.ActiveDocument\n\n$needle = \"al igual que los clientes que ganaron en dicho sorteo\"\n$insertion = \", en total hay 3 sorteos diarios (esto puede variar)\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n\nif (-not $found) {\n    throw \"Target sentence not found in document.\"\n}\n\n# $rng now spans exactly the matched text (no trailing period). Collapse to\n# its end point (wdCollapseEnd = 0) so the caret sits right before the \".\".\n$rng.Collapse(0)\n$rng.InsertBefore($insertion)\n\n$d.Save()\n"}
